# Idempotent script to ensure the payroll import template contains the
# standard sample rows:
#   - Models sheet : a weekly-paid sample model "M-003"
#   - Payouts sheet: an on_hold sample payout for "M-003"
#
# Running this more than once must not create duplicate rows.
#
# NOTE: reading via the ".Value" getter is unreliable in this runtime, so
# all reads below use ".Value2" instead; ".Value" is still used (and
# works correctly) for writes.

$wb = $excel.ActiveWorkbook

$sampleCode = "M-003"

# --- helper: find the last used row in column A of a worksheet --------
function Get-LastRow($sheet) {
    if ($sheet.Cells.Item(1, 1).Value2 -eq $null) {
        return 0
    }
    $lastCell = $sheet.Cells.Item($sheet.Rows.Count, 1).End(-4162)
    return $lastCell.Row
}

# --- Models sheet: add the sample weekly model if it is missing -------
$modelsSheet = $wb.Worksheets.Item("Models")
$modelsLastRow = Get-LastRow $modelsSheet

$modelsHasSample = $false
for ($r = 2; $r -le $modelsLastRow; $r++) {
    if ($modelsSheet.Cells.Item($r, 1).Value2 -eq $sampleCode) {
        $modelsHasSample = $true
        break
    }
}

if (-not $modelsHasSample) {
    $newRow = $modelsLastRow + 1

    $modelsSheet.Cells.Item($newRow, 1).Value = "M-003"
    $modelsSheet.Cells.Item($newRow, 2).Value = "Active"
    $modelsSheet.Cells.Item($newRow, 3).Value = "Sample Weekly Model"
    $modelsSheet.Cells.Item($newRow, 4).Value = "WeeklySample"

    # Keep the date as literal text (matches existing "YYYY-MM-DD" sample
    # cells) instead of letting Excel coerce it into a date serial value.
    $modelsSheet.Cells.Item($newRow, 5).NumberFormat = "@"
    $modelsSheet.Cells.Item($newRow, 5).Value = "2025-11-04"
    $modelsSheet.Cells.Item($newRow, 5).Style = "Normal"

    $modelsSheet.Cells.Item($newRow, 6).Value = "Bank Transfer"
    $modelsSheet.Cells.Item($newRow, 7).Value = "weekly"
    $modelsSheet.Cells.Item($newRow, 8).Value = 1000
}

# --- Payouts sheet: add the sample on_hold payout if it is missing ----
$payoutsSheet = $wb.Worksheets.Item("Payouts")
$payoutsLastRow = Get-LastRow $payoutsSheet

$payoutsHasSample = $false
for ($r = 2; $r -le $payoutsLastRow; $r++) {
    $codeMatches = $payoutsSheet.Cells.Item($r, 1).Value2 -eq $sampleCode
    $statusMatches = $payoutsSheet.Cells.Item($r, 4).Value2 -eq "on_hold"
    if ($codeMatches -and $statusMatches) {
        $payoutsHasSample = $true
        break
    }
}

if (-not $payoutsHasSample) {
    $newRow = $payoutsLastRow + 1

    $payoutsSheet.Cells.Item($newRow, 1).Value = "M-003"

    $payoutsSheet.Cells.Item($newRow, 2).NumberFormat = "@"
    $payoutsSheet.Cells.Item($newRow, 2).Value = "2025-11-04"
    $payoutsSheet.Cells.Item($newRow, 2).Style = "Normal"

    $payoutsSheet.Cells.Item($newRow, 3).Value = 250
    $payoutsSheet.Cells.Item($newRow, 4).Value = "on_hold"
    $payoutsSheet.Cells.Item($newRow, 5).Value = "Bank Transfer"
}
